# Apply edit: add columns I (I0) and J (IF) with per-row values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): I1 = "I0", J1 = "IF", matching the bold/bordered
# header style already used by B1:H1 (style index 1), copied from H1 ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2..67: column I (I0) and column J (IF) values ---
$iValues = @(6,8,7,7,7,9,8,8,7,7,7,7,8,8,7,8,9,6,6,7,8,7,7,7,7,8,8,7,7,6,8,9,8,8,8,7,8,10,8,7,8,9,6,6,5,7,8,9,7,9,8,9,9,7,8,8,7,7,7,8,6,4,8,6,7,5)
$jValues = @(6,8,7,7,7,9,8,8,7,7,7,7,8,8,7,8,9,6,6,7,8,7,7,7,7,8,8,7,7,7,8,9,8,8,8,7,9,10,8,8,8,9,6,6,6,8,8,9,7,9,8,9,9,7,8,8,7,7,7,8,6,4,8,6,7,5)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}

Write-Host "Applied I0/IF columns to rows 1-67"
